# Poland IV Liga - "Atualizacao de bases das ligas" update
#
# The underlying dataset had a handful of match rows that needed to trade
# places with their neighbour (the row's rank/Id in column A and the match
# Date in column D stay put, but every other column - match Id, HomeTeam,
# AwayTeam, scores and all the odds columns - belongs to the other row).
#
# For each pair below we swap the full B:AD block between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows {
    param($RowA, $RowB)

    $rangeA = $ws.Range("B$($RowA):AD$($RowA)")
    $rangeB = $ws.Range("B$($RowB):AD$($RowB)")

    $valuesA = $rangeA.Value()
    $valuesB = $rangeB.Value()

    $rangeA.Value = $valuesB
    $rangeB.Value = $valuesA
}

Swap-Rows 28 29
Swap-Rows 75 76
Swap-Rows 131 132
Swap-Rows 215 216
Swap-Rows 222 223
